# Version 2.0.1 - solucionado error espera de base de datos
# Updates the patient admission/discharge form (HOJA DE INGRESO Y EGRESO)
# with the new patient's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Patient identification (row 6) ---
$ws.Range("A6").Value = "Pérez"
$ws.Range("C6").Value = "Hernández"
$ws.Range("E6").Value = "Anthony"
$ws.Range("G6").Value = "Alejandro"
$ws.Range("I6").Value = "/201761947"

# --- Dirección actual (row 8) ---
$ws.Range("A8").Value = "28 ave. 23-69"
$ws.Range("D8").Value = "Z. 7 Col. 4 de febrero"
$ws.Range("F8").Value = "Guatemal"
$ws.Range("H8").Value = "Guatemala"
$ws.Range("J8").Value = ""

# --- Fecha de nacimiento / Lugar de nacimiento / Sexo (row 12) ---
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "2014-10-29"
$ws.Range("H12").Value = "Guatemala"
$ws.Range("J12").Value = "Masculino"

# --- Nacionalidad / No. de Cédula (row 14) ---
$ws.Range("F14").Value = "Guatemalteco"
$ws.Range("H14").Value = "Ins. Nac. 331027"

# --- Nombre del Padre / Nombre de la Madre (row 18) ---
$ws.Range("A18").Value = "Fredy Misael Pérez García"
$ws.Range("F18").Value = "Miryam Angélica Hernández González"

# --- En caso de emergencia notificar a (row 20) ---
$ws.Range("A20").Value = "Miryam Hernández"
$ws.Range("F20").Value = "Madre"
$ws.Range("J20").NumberFormat = "@"
$ws.Range("J20").Value = "33315313"

# --- Hora de ingreso (row 24) ---
$ws.Range("C24").Value = "12:20:11"
